# PacketRule.xlsx bug-fix edit
# - fix typo in sheet name: "ClilentPacketFactory" -> "ClientPacketFactory"
# - move the active selection on "PacketData" from B19 to D19
# - move the active tab from "PacketClass" to "ClientPacketFactory"

$wb = $excel.ActiveWorkbook

# Fix the misspelled sheet name
$wsFactory = $wb.Worksheets.Item("ClilentPacketFactory")
$wsFactory.Name = "ClientPacketFactory"

# Update the stored selection on PacketData (B19 -> D19)
$wsData = $wb.Worksheets.Item("PacketData")
$null = $wsData.Activate()
$null = $wsData.Range("D19").Select()

# Finally, activate ClientPacketFactory so it becomes the active/selected tab
$null = $wsFactory.Activate()
